$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Con"
$ws.Range("C5").Value = "Cat"
$ws.Range("D5").Formula = "=CONCAT(B5:C5)"

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

[void]$ws.Range("D5").Select()
